# Auto-generated Excel COM-interop script
# Applies the numeric corrections described in the commit diff
# to worksheets ALC, ARM, CRP, CUL, GSM, LTW, WVR.

$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 192.3077
$ws.Range("I9").Value = 192.1
$ws.Range("K9").Value = 192.1
$ws.Range("M9").Value = -23.09999999999999
$ws.Range("H38").Value = 305.42856
$ws.Range("I38").Value = 366.9091
$ws.Range("J38").Value = 80
$ws.Range("K38").Value = 1100.7273
$ws.Range("L38").Value = 240
$ws.Range("M38").Value = -728.7273
$ws.Range("N38").Value = -984
$ws.Range("H42").Value = 269.2
$ws.Range("J42").Value = 390.33334
$ws.Range("L42").Value = 1171.00002
$ws.Range("N42").Value = -1631.00002
$ws.Range("H87").Value = 33348.6
$ws.Range("J87").Value = 33348.6
$ws.Range("L87").Value = 33348.6
$ws.Range("N87").Value = -35844.6
$ws.Range("H90").Value = 33348.6
$ws.Range("J90").Value = 33348.6
$ws.Range("L90").Value = 100045.8
$ws.Range("N90").Value = -112525.8
$ws.Range("H138").Value = 1838.6111
$ws.Range("I138").Value = 1393.5
$ws.Range("J138").Value = 5399.5
$ws.Range("K138").Value = 4180.5
$ws.Range("L138").Value = 16198.5
$ws.Range("M138").Value = 959.5
$ws.Range("N138").Value = -26478.5

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 4000
$ws.Range("I25").Value = 4000
$ws.Range("J25").Value = 4000
$ws.Range("K25").Value = 4000
$ws.Range("L25").Value = 4000
$ws.Range("M25").Value = -3598
$ws.Range("N25").Value = -4804
$ws.Range("H32").Value = 46041.21
$ws.Range("I32").Value = 27712
$ws.Range("K32").Value = 27712
$ws.Range("M32").Value = -27425
$ws.Range("H61").Value = 5000
$ws.Range("I61").Value = 5000
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 5000
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4788
$ws.Range("N61").ClearContents()
$ws.Range("H132").Value = 20483.084
$ws.Range("I132").Value = 24977.445
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 74932.33499999999
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -72402.33499999999
$ws.Range("N132").Value = -26060
$ws.Range("H136").Value = 5000
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -12450
$ws.Range("N136").ClearContents()

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 5624.7827
$ws.Range("I58").Value = 7129.5
$ws.Range("K58").Value = 7129.5
$ws.Range("M58").Value = -6926.5
$ws.Range("H132").Value = 2556.919
$ws.Range("I132").Value = 2229.12
$ws.Range("K132").Value = 6687.36
$ws.Range("M132").Value = -4157.36
$ws.Range("H134").Value = 2211.7693
$ws.Range("I134").Value = 1979.4166
$ws.Range("K134").Value = 5938.2498
$ws.Range("M134").Value = -3403.2498
$ws.Range("H136").Value = 5624.7827
$ws.Range("I136").Value = 7129.5
$ws.Range("K136").Value = 21388.5
$ws.Range("M136").Value = -18838.5

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 9499.5
$ws.Range("J63").Value = 14999
$ws.Range("L63").Value = 44997
$ws.Range("N63").Value = -46495
$ws.Range("H66").Value = 9499.5
$ws.Range("J66").Value = 14999
$ws.Range("L66").Value = 134991
$ws.Range("N66").Value = -142479

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 23271
$ws.Range("J64").Value = 23271
$ws.Range("L64").Value = 23271
$ws.Range("N64").Value = -23767
$ws.Range("H67").Value = 23271
$ws.Range("J67").Value = 23271
$ws.Range("L67").Value = 23271
$ws.Range("N67").Value = -24987
$ws.Range("H70").Value = 14290404
$ws.Range("I70").Value = 16671388
$ws.Range("J70").Value = 4500
$ws.Range("K70").Value = 16671388
$ws.Range("L70").Value = 4500
$ws.Range("M70").Value = -16671118
$ws.Range("N70").Value = -5040
$ws.Range("H73").Value = 14290404
$ws.Range("I73").Value = 16671388
$ws.Range("J73").Value = 4500
$ws.Range("K73").Value = 16671388
$ws.Range("L73").Value = 4500
$ws.Range("M73").Value = -16670452
$ws.Range("N73").Value = -6372
$ws.Range("H132").Value = 6196.8335
$ws.Range("I132").Value = 8795.333000000001
$ws.Range("K132").Value = 26385.999
$ws.Range("M132").Value = -23855.999

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value = 13214.429
$ws.Range("I7").Value = 15168
$ws.Range("J7").Value = 11749.25
$ws.Range("K7").Value = 15168
$ws.Range("L7").Value = 11749.25
$ws.Range("M7").Value = -15056
$ws.Range("N7").Value = -11973.25
$ws.Range("H24").Value = 5000
$ws.Range("J24").Value = 5000
$ws.Range("L24").Value = 5000
$ws.Range("N24").Value = -5686
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H40").Value = 6994.6665
$ws.Range("I40").Value = 4995
$ws.Range("K40").Value = 4995
$ws.Range("M40").Value = -4859
$ws.Range("H122").Value = 4386.3335
$ws.Range("I122").Value = 3996.5
$ws.Range("K122").Value = 11989.5
$ws.Range("M122").Value = -9539.5
$ws.Range("H126").Value = 13214.429
$ws.Range("I126").Value = 15168
$ws.Range("J126").Value = 11749.25
$ws.Range("K126").Value = 45504
$ws.Range("L126").Value = 35247.75
$ws.Range("M126").Value = -43034
$ws.Range("N126").Value = -40187.75
$ws.Range("H136").Value = 3606.5833
$ws.Range("I136").Value = 3098
$ws.Range("K136").Value = 9294
$ws.Range("M136").Value = -6744
$ws.Range("H140").Value = 94250
$ws.Range("J140").Value = 94250
$ws.Range("L140").Value = 94250
$ws.Range("N140").Value = -104610

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 11666.667
$ws.Range("I29").Value = 11666.667
$ws.Range("K29").Value = 11666.667
$ws.Range("M29").Value = -11376.667
$ws.Range("H122").Value = 1868.9
$ws.Range("I122").Value = 906.2857
$ws.Range("J122").Value = 2387.2307
$ws.Range("K122").Value = 2718.8571
$ws.Range("L122").Value = 7161.6921
$ws.Range("M122").Value = -268.8571000000002
$ws.Range("N122").Value = -12061.6921
$ws.Range("H126").Value = 3605.4285
$ws.Range("J126").Value = 3797.8
$ws.Range("L126").Value = 11393.4
$ws.Range("N126").Value = -16333.4
$ws.Range("H132").Value = 14912.25
$ws.Range("I132").Value = 17976.4
$ws.Range("K132").Value = 53929.2
$ws.Range("M132").Value = -51399.2
